$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.075.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.86%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.105.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.55%  '

$ws.Range("E6").Value = '  -0.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5180'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4467'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09521'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.177'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.109.68'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.753'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.095'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001173'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.009'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06705'
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = '  -0.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.190'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.153.85'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.323'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.358.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.536'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.160'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1056'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.630'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.256'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.945'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.179'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.17'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02573'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06783'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2284'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6958'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.310'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.95%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6708'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.285'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.34%  '

$ws.Range("E47").Value = '  -1.85%  '

$ws.Range("E48").Value = '  -3.47%  '

$ws.Range("E49").Value = '  -2.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07186'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.71%  '
